$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sale (row 122, Dandenong Station bay "10-11") now notes it is the
# TRBS and coach stop.
$ws.Range("D122").Value = "TRBS and coach stop"

# Knox City SC/Burwood Hwy (Wantirna South) bays (rows 210-219):
# the old numeric "bay order" values (1-10) are replaced with their
# current bay letters, and the former numeric bay labels are preserved
# in a new column D ("Former Bay N").
$ws.Range("B210").Value = "D"
$ws.Range("D210").Value = "Former Bay 1"

$ws.Range("B211").Value = "E"
$ws.Range("D211").Value = "Former Bay 2"

$ws.Range("B212").Value = "F"
$ws.Range("D212").Value = "Former Bay 3"

$ws.Range("B213").Value = "G"
$ws.Range("D213").Value = "Former Bay 4"

$ws.Range("B214").Value = "H"
$ws.Range("D214").Value = "Former Bay 5"

$ws.Range("B215").Value = "A"
$ws.Range("D215").Value = "Former Bay 6"

$ws.Range("B216").Value = "B"
$ws.Range("D216").Value = "Former Bay 7"

$ws.Range("B217").Value = "C"
$ws.Range("D217").Value = "Former Bay 8"

$ws.Range("B218").Value = "J"
$ws.Range("D218").Value = "Former Bay 9"

$ws.Range("B219").Value = "I"
$ws.Range("D219").Value = "Former Bay 10"

# Restore the view to where the author left off editing.
$ws.Application.ActiveWindow.ScrollRow = 197
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.Zoom = 141
$ws.Range("B219").Select()
